# Updated cryptos list on Thu Jan 18 09:30:22 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on Sheet1, and fixes the Chainlink / WrappedEther rows which
# had swapped places in the source feed (rows 15 and 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain text (values like "42.721.88" use '.' as a
# thousands separator as well as a decimal point, so they can never be real
# numbers) -- but plenty of the individual numbers we're writing here (e.g.
# "312.04") WOULD parse as a normal float. Excel's Range.Value setter
# auto-detects that and silently turns the cell into a Number, which would
# not match the original text-cell layout. Forcing the format to Text first
# (and resetting the style afterwards so we don't leave a stray style on the
# cell) keeps these as literal text, exactly like the rest of the column.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row 2: Bitcoin ---
Set-TextValue 'D2' '42.785.05'
$ws.Range('E2').Value = '  -0.02%  '

# --- Row 3: Ethereum ---
Set-TextValue 'D3' '2.543.93'
$ws.Range('E3').Value = '  -0.19%  '

# --- Row 4: TetherUSD ---
$ws.Range('E4').Value = '  -0.01%  '

# --- Row 5: BNB ---
Set-TextValue 'D5' '312.04'
$ws.Range('E5').Value = '  +0.58%  '

# --- Row 6: Solana ---
Set-TextValue 'D6' '100.62'
$ws.Range('E6').Value = '  +2.03%  '

# --- Row 7: XRP ---
$ws.Range('E7').Value = '  -0.82%  '

# --- Row 8: USDC ---
$ws.Range('E8').Value = '  +0.02%  '

# --- Row 9: Cardano ---
$ws.Range('E9').Value = '  -1.70%  '

# --- Row 10: Avalanche ---
Set-TextValue 'D10' '35.51'
$ws.Range('E10').Value = '  -0.89%  '

# --- Row 11: Dogecoin ---
$ws.Range('E11').Value = '  -0.37%  '

# --- Row 12: Polkadot ---
$ws.Range('E12').Value = '  -1.43%  '

# --- Row 13: TRON ---
$ws.Range('E13').Value = '  +0.98%  '

# --- Row 14: WrappedliquidstakedEther2.0 ---
Set-TextValue 'D14' '2.930.57'
$ws.Range('E14').Value = '  -0.33%  '

# --- Row 15 & 16: Chainlink and WrappedEther swapped positions ---
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '15.40'
$ws.Range('E15').Value = '  -3.45%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '2.542.67'
$ws.Range('E16').Value = '  -0.16%  '

# --- Row 17: Polygon ---
Set-TextValue 'D17' '0.818'
$ws.Range('E17').Value = '  -2.73%  '

# --- Row 18: WrappedBTC ---
Set-TextValue 'D18' '42.781.71'
$ws.Range('E18').Value = '  -0.12%  '

# --- Row 19: Uniswap ---
Set-TextValue 'D19' '6.74'
$ws.Range('E19').Value = '  -0.20%  '

# --- Row 20: InternetComputer(DFINITY) ---
Set-TextValue 'D20' '12.33'
$ws.Range('E20').Value = '  -0.91%  '

# --- Row 21: ShibaInu ---
$ws.Range('E21').Value = '  -0.49%  '

# --- Row 22: Litecoin ---
Set-TextValue 'D22' '69.95'
$ws.Range('E22').Value = '  +0.75%  '

# --- Row 23: BitcoinCash ---
Set-TextValue 'D23' '243.21'
$ws.Range('E23').Value = '  -2.17%  '

# --- Row 24: PancakeSwap ---
$ws.Range('E24').Value = '  -1.52%  '

# --- Row 25: ImmutableX ---
$ws.Range('E25').Value = '  -1.68%  '

# --- Row 26: Dai (unchanged) ---

# --- Row 27: EthereumClassic ---
Set-TextValue 'D27' '25.66'
$ws.Range('E27').Value = '  -3.97%  '

# --- Row 28: Toncoin ---
$ws.Range('E28').Value = '  -1.22%  '

# --- Row 29: Cosmos ---
$ws.Range('E29').Value = '  -0.14%  '

# --- Row 30: InjectiveProtocol ---
Set-TextValue 'D30' '38.43'
$ws.Range('E30').Value = '  -4.34%  '

# --- Row 31: Filecoin ---
Set-TextValue 'D31' '5.88'
$ws.Range('E31').Value = '  +2.00%  '

# --- Row 32: Monero ---
Set-TextValue 'D32' '157.97'
$ws.Range('E32').Value = '  -0.44%  '

# --- Row 33: ApeXProtocol ---
$ws.Range('E33').Value = '  +5.78%  '

# --- Row 34: WEMIXToken ---
$ws.Range('E34').Value = '  +1.87%  '

# --- Row 35: Hedera ---
Set-TextValue 'D35' '0.0794'
$ws.Range('E35').Value = '  -0.98%  '

# --- Row 36: LidoDAOToken ---
$ws.Range('E36').Value = '  -4.17%  '

# --- Row 37: Celestia ---
Set-TextValue 'D37' '17.98'
$ws.Range('E37').Value = '  -2.02%  '

# --- Row 38: ARBITRUM ---
$ws.Range('E38').Value = '  -5.45%  '

# --- Row 39: Kaspa ---
$ws.Range('E39').Value = '  -0.51%  '

# --- Row 40: Stellar ---
$ws.Range('E40').Value = '  -0.34%  '

# --- Row 41: RenderToken ---
Set-TextValue 'D41' '4.14'
$ws.Range('E41').Value = '  +0.62%  '

# --- Row 42: EnergySwap ---
Set-TextValue 'D42' '22.01'
$ws.Range('E42').Value = '  -2.44%  '

# --- Row 43: FirstDigitalUSD ---
$ws.Range('E43').Value = '  +0.12%  '

# --- Row 44: VeChain ---
$ws.Range('E44').Value = '  -0.49%  '

# --- Row 45: NEARProtocol ---
$ws.Range('E45').Value = '  +2.03%  '

# --- Row 46: Maker ---
Set-TextValue 'D46' '1.998.60'
$ws.Range('E46').Value = '  +0.06%  '

# --- Row 47: FraxShare ---
Set-TextValue 'D47' '9.15'
$ws.Range('E47').Value = '  +0.51%  '

# --- Row 48: RocketPoolETH ---
Set-TextValue 'D48' '2.788.23'
$ws.Range('E48').Value = '  -0.02%  '

# --- Row 49: Algorand ---
$ws.Range('E49').Value = '  -0.16%  '

# --- Row 50: BitcoinSV ---
$ws.Range('E50').Value = '  -1.79%  '

# --- Row 51: ordi ---
Set-TextValue 'D51' '72.48'
$ws.Range('E51').Value = '  -1.60%  '
